$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 169
$ws.Cells.Item(6, 6).Value = 345
$ws.Cells.Item(7, 6).Value = 5981
$ws.Cells.Item(8, 6).Value = 10188
$ws.Cells.Item(9, 6).Value = 4028
$ws.Cells.Item(10, 6).Value = 222
$ws.Cells.Item(11, 6).Value = 146
$ws.Cells.Item(12, 6).Value = 62
$ws.Cells.Item(13, 6).Value = 137
$ws.Cells.Item(15, 6).Value = 4021
$ws.Cells.Item(17, 6).Value = 152
$ws.Cells.Item(19, 6).Value = 5663
$ws.Cells.Item(21, 6).Value = 2222
$ws.Cells.Item(22, 6).Value = 149
$ws.Cells.Item(23, 6).Value = 407
$ws.Cells.Item(24, 6).Value = 8497
$ws.Cells.Item(26, 6).Value = 90
$ws.Cells.Item(27, 6).Value = 2249
$ws.Cells.Item(28, 6).Value = 2294
$ws.Cells.Item(29, 6).Value = 1356
$ws.Cells.Item(31, 6).Value = 1873
$ws.Cells.Item(32, 6).Value = 41
$ws.Cells.Item(33, 6).Value = 300
$ws.Cells.Item(35, 6).Value = 19
$ws.Cells.Item(36, 6).Value = 273
$ws.Cells.Item(37, 6).Value = 31
$ws.Cells.Item(38, 6).Value = 21
$ws.Cells.Item(41, 6).Value = 55
$ws.Cells.Item(42, 6).Value = 77
$ws.Cells.Item(43, 6).Value = 203
$ws.Cells.Item(44, 6).Value = 1409
$ws.Cells.Item(45, 6).Value = 2278
$ws.Cells.Item(46, 6).Value = 168
$ws.Cells.Item(47, 6).Value = 253
$ws.Cells.Item(48, 6).Value = 1234
$ws.Cells.Item(49, 6).Value = 14

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(5, 6).Value = 156

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 637

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 169
$ws.Cells.Item(4, 6).Value = 637
$ws.Cells.Item(6, 6).Value = 345
$ws.Cells.Item(7, 6).Value = 5981
$ws.Cells.Item(8, 6).Value = 4028
$ws.Cells.Item(9, 6).Value = 222
$ws.Cells.Item(10, 6).Value = 146
$ws.Cells.Item(11, 6).Value = 62
$ws.Cells.Item(12, 6).Value = 137
$ws.Cells.Item(13, 6).Value = 156
$ws.Cells.Item(16, 6).Value = 4021
$ws.Cells.Item(18, 6).Value = 152
$ws.Cells.Item(20, 6).Value = 5663
$ws.Cells.Item(22, 6).Value = 2222
$ws.Cells.Item(23, 6).Value = 149
$ws.Cells.Item(24, 6).Value = 407
$ws.Cells.Item(25, 6).Value = 8497
$ws.Cells.Item(28, 6).Value = 2249
$ws.Cells.Item(29, 6).Value = 2294
$ws.Cells.Item(30, 6).Value = 1356
$ws.Cells.Item(32, 6).Value = 1873
$ws.Cells.Item(33, 6).Value = 41
$ws.Cells.Item(34, 6).Value = 300
$ws.Cells.Item(35, 6).Value = 19
$ws.Cells.Item(36, 6).Value = 273
$ws.Cells.Item(37, 6).Value = 31
$ws.Cells.Item(38, 6).Value = 21
$ws.Cells.Item(41, 6).Value = 55
$ws.Cells.Item(42, 6).Value = 77
$ws.Cells.Item(43, 6).Value = 203
$ws.Cells.Item(44, 6).Value = 1409
$ws.Cells.Item(46, 6).Value = 2278
$ws.Cells.Item(47, 6).Value = 168
$ws.Cells.Item(48, 6).Value = 253
$ws.Cells.Item(49, 6).Value = 1234
